# Apply the 2026-03-02 snapshot refresh to the Betfair Back/Lay sheet.
#
# Semantics, derived from the OOXML diff:
#  - Row 14 (Spanish La Liga / Real Madrid v Getafe) is gone; every row
#    below it shifts up by one, and the trailing row (old row 19) disappears
#    -> net effect: delete row 14 entirely (Excel's own row-delete already
#    renumbers everything and shrinks the used range to A1:BH18).
#  - A handful of odds cells get refreshed values on top of that shift.
#  - Every remaining row's SnapshotTS (column BH) moves to the later time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the Real Madrid v Getafe row; everything below shifts up and the
#    dimension shrinks from A1:BH19 to A1:BH18 automatically.
$ws.Rows.Item(14).Delete()

# 2) Odds refreshed on rows that did NOT move (still 2-13 after the delete).
$ws.Range("F2").Value = 1.35
$ws.Range("K2").Value = 8.6

$ws.Range("J4").Value = 2.98

$ws.Range("F5").Value = 4.3
$ws.Range("G5").Value = 4.5
$ws.Range("K5").Value = 3.5

$ws.Range("F13").Value = 3.25
$ws.Range("G13").Value = 3.35
$ws.Range("AQ13").Value = 9
$ws.Range("AU13").Value = 7
$ws.Range("BB13").Value = 29

# 3) Odds refreshed on rows that now sit one row higher than before
#    (old 15->14, old 16->15, old 17->16, old 18->17, old 19->18).
# New row 14 (Gil Vicente v Benfica)
$ws.Range("F14").Value = 2.48
$ws.Range("H14").Value = 1.61
$ws.Range("K14").Value = 950
$ws.Range("P14").Value = 1.71
$ws.Range("Q14").Value = 1.92

# New row 16 (Estudiantes v Velez Sarsfield)
$ws.Range("P16").Value = 1.42

# New row 17 (Deportivo Riestra v CA Platense)
$ws.Range("P17").Value = 1.28
$ws.Range("Q17").Value = 3.9

# 4) Refresh SnapshotTS on every remaining row (2-18).
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 60).Value = "2026-02-28 07:40:36"
}
